$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (QLDGovernment) - add Lab/Morph/Sites = OK
$ws.Range("B10").Value = "OK"
$ws.Range("C10").Value = "OK"
$ws.Range("D10").Value = "OK"

# Row 11 (NatGeoChemicalSurvey) - add Lab/Morph/Sites = OK, remove old comment in E11
$ws.Range("B11").Value = "OK"
$ws.Range("C11").Value = "OK"
$ws.Range("D11").Value = "OK"
$ws.Range("E11").ClearContents()

# New dataset rows - enter column A labels first (matches the order the
# new shared strings were introduced in the authored workbook)
$ws.Range("A12").Value = "TasGovernment"
$ws.Range("A13").Value = "SAGovernment"
$ws.Range("A14").Value = "EcologicalProjects"
$ws.Range("A15").Value = "EastCentral_Australia"

# Row 12 (TasGovernment) - Lab/Morph/Sites = nssc
$ws.Range("B12").Value = "nssc"
$ws.Range("C12").Value = "nssc"
$ws.Range("D12").Value = "nssc"

# Row 13 (SAGovernment) - Lab/Morph/Sites = nssc
$ws.Range("B13").Value = "nssc"
$ws.Range("C13").Value = "nssc"
$ws.Range("D13").Value = "nssc"

# Row 14 (EcologicalProjects) has no further data yet

# Row 15 (EastCentral_Australia) - Lab/Morph/Sites = OK
$ws.Range("B15").Value = "OK"
$ws.Range("C15").Value = "OK"
$ws.Range("D15").Value = "OK"

# Update selection to match final state
$null = $ws.Range("A14").Select()
